# Insert a new "ListBullet" paragraph with the new docente's name right
# after the "Docente(s) Responsável(eis)" heading paragraph.

$d = $word.ActiveDocument

# Locate the "Docente(s) Responsável(eis)" heading paragraph.
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "Docente(s) Responsável(eis)*") {
        $target = $p
    }
}

# Collapse to the end of that paragraph and insert a brand-new paragraph
# right after it.
$r = $target.Range
$r.Collapse(0)
$r.InsertParagraphAfter()

# The newly inserted paragraph is now the very next paragraph after the
# heading; fill in its text and apply the bullet-list style used by the
# other "Docente(s)" style entries in this document.
$newPara = $target.Next()
$newPara.Range.Text = "5464150 - Mariana Consiglio Kasemodel"
$newPara.Style = "ListBullet"
